# Apply "Natmi following Dr Hou advice" update: rows 2-6 revised, rows 7-11 added
# (full cross-product of Sending cluster x Target cluster in {ECs, M2})
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A="ECs"; D="ECs"; E=3; F=1; G=43.64525366666666; H=130.935761; I=0.5421200012818527; J=0.5463095076486307; K=3; L=1; M=9.112632333333332; N=27.337897; O=0.9981738658344552; P=0.9981738658344552; Q=397.7231497594018; R=3579.508347834616; S=0.5411300174256867; T=0.5453118731917517 },
    @{ Row=3; A="ECs"; D="M2"; E=3; F=1; G=43.64525366666666; H=130.935761; I=0.5421200012818527; J=0.5463095076486307; K=1; L=0.3333333333333333; M=0.01667133333333333; N=0.050014; O=0.001826134165544791; P=0.001826134165544791; Q=0.7276245722948887; R=6.548621150653999; S=0.0009899838561659771; T=0.0009976344568791179 },
    @{ Row=4; A="FAPs"; D="ECs"; E=3; F=1; G=0.161567; H=0.484701; I=0.00200683224150899; J=0.00202234105216526; K=3; L=1; M=9.112632333333332; N=27.337897; O=0.9981738658344552; P=0.9981738658344552; Q=1.472300668199667; R=13.250706013797; S=0.002003167496588253; T=0.002018647986075517 },
    @{ Row=5; A="FAPs"; D="M2"; E=3; F=1; G=0.161567; H=0.484701; I=0.00200683224150899; J=0.00202234105216526; K=1; L=0.3333333333333333; M=0.01667133333333333; N=0.050014; O=0.001826134165544791; P=0.001826134165544791; Q=0.002693537312666667; R=0.024241835814; S=0.000003664744920736402; T=0.000003693066089742782 },
    @{ Row=6; A="M1"; D="ECs"; E=3; F=1; G=10.44622; H=31.33866; I=0.1297530504242576; J=0.1307557827152189; K=3; L=1; M=9.112632333333332; N=27.337897; O=0.9981738658344552; P=0.9981738658344552; Q=95.19256213311331; R=856.7330591980199; S=0.1295161039457942; T=0.1305170051130601 },
    @{ Row=7; A="M1"; D="M2"; E=3; F=1; G=10.44622; H=31.33866; I=0.1297530504242576; J=0.1307557827152189; K=1; L=0.3333333333333333; M=0.01667133333333333; N=0.050014; O=0.001826134165544791; P=0.001826134165544791; Q=0.1741524156933333; R=1.56737174124; S=0.0002369464784633929; T=0.0002387776021588123 },
    @{ Row=8; A="M2"; D="ECs"; E=3; F=1; G=24.40323666666667; H=73.20971; I=0.303113891697197; J=0.3054563575278647; K=3; L=1; M=9.112632333333332; N=27.337897; O=0.9981738658344552; P=0.9981738658344552; Q=222.3777234866522; R=2001.39951137987; S=0.3025603650635175; T=0.3048985532373002 },
    @{ Row=9; A="M2"; D="M2"; E=3; F=1; G=24.40323666666667; H=73.20971; I=0.303113891697197; J=0.3054563575278647; K=1; L=0.3333333333333333; M=0.01667133333333333; N=0.050014; O=0.001826134165544791; P=0.001826134165544791; Q=0.4068344928822222; R=3.66151043594; S=0.0005535266336794951; T=0.0005578042905644985 },
    @{ Row=10; A="sCs"; D="ECs"; E=2; F=1; G=1.852196; H=3.704392; I=0.02300622435518382; J=0.01545601105612031; K=3; L=1; M=9.112632333333332; N=27.337897; O=0.9981738658344552; P=0.9981738658344552; Q=16.87838115727067; R=101.270286943624; S=0.02296421190286863; T=0.01542778630626769 },
    @{ Row=11; A="sCs"; D="M2"; E=2; F=1; G=1.852196; H=3.704392; I=0.02300622435518382; J=0.01545601105612031; K=1; L=0.3333333333333333; M=0.01667133333333333; N=0.050014; O=0.001826134165544791; P=0.001826134165544791; Q=0.03087857691466667; R=0.185271461488; S=0.00004201245231518986; T=0.00002822474985261933 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("A$i").Value = $r.A
    $ws.Range("B$i").Value = "Ceacam1"
    $ws.Range("C$i").Value = "Sele"
    $ws.Range("D$i").Value = $r.D
    $ws.Range("E$i").Value = $r.E
    $ws.Range("F$i").Value = $r.F
    $ws.Range("G$i").Value = $r.G
    $ws.Range("H$i").Value = $r.H
    $ws.Range("I$i").Value = $r.I
    $ws.Range("J$i").Value = $r.J
    $ws.Range("K$i").Value = $r.K
    $ws.Range("L$i").Value = $r.L
    $ws.Range("M$i").Value = $r.M
    $ws.Range("N$i").Value = $r.N
    $ws.Range("O$i").Value = $r.O
    $ws.Range("P$i").Value = $r.P
    $ws.Range("Q$i").Value = $r.Q
    $ws.Range("R$i").Value = $r.R
    $ws.Range("S$i").Value = $r.S
    $ws.Range("T$i").Value = $r.T
}
